$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3479973333333333
$ws.Range("H2").Value = 1.043992
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.920975666666667
$ws.Range("N2").Value = 29.762927
$ws.Range("O2").Value = 0.5379109985456686
$ws.Range("P2").Value = 0.5379109985456685
$ws.Range("Q2").Value = 3.452473076064889
$ws.Range("R2").Value = 31.072257684584
$ws.Range("S2").Value = 0.5379109985456686
$ws.Range("T2").Value = 0.5379109985456685

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3479973333333333
$ws.Range("H3").Value = 1.043992
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.153936333333334
$ws.Range("N3").Value = 18.461809
$ws.Range("O3").Value = 0.3336637594195427
$ws.Range("P3").Value = 0.3336637594195427
$ws.Range("Q3").Value = 2.141553433503111
$ws.Range("R3").Value = 19.273980901528
$ws.Range("S3").Value = 0.3336637594195427
$ws.Range("T3").Value = 0.3336637594195427

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3479973333333333
$ws.Range("H4").Value = 1.043992
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.368614333333333
$ws.Range("N4").Value = 7.105843
$ws.Range("O4").Value = 0.1284252420347888
$ws.Range("P4").Value = 0.1284252420347888
$ws.Range("Q4").Value = 0.8242714716951109
$ws.Range("R4").Value = 7.418443245255999
$ws.Range("S4").Value = 0.1284252420347888
$ws.Range("T4").Value = 0.1284252420347888
